# Daily refresh of the cryptos price/volume table (GitHub Actions job).
# Column layout: A=rank(idx) B=Coin C=Link D=Price E=Volume(1h)
# For D-column cells whose new text parses as a plain number, force the
# cell to Text format first so Excel keeps the original string (e.g.
# "1.00", "597.53") instead of collapsing it to a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.003.18'
$ws.Cells.Item(2, 5).Value = '  -1.19%  '
$ws.Cells.Item(3, 4).Value = '3.847.63'
$ws.Cells.Item(3, 5).Value = '  -1.68%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.12%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '597.53'
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '166.96'
$ws.Cells.Item(6, 5).Value = '  +0.85%  '
$ws.Cells.Item(7, 4).Value = '3.847.47'
$ws.Cells.Item(7, 5).Value = '  -1.72%  '
$ws.Cells.Item(8, 5).Value = '  +0.16%  '
$ws.Cells.Item(9, 5).Value = '  -0.91%  '
$ws.Cells.Item(10, 5).Value = '  -0.78%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.30'
$ws.Cells.Item(11, 5).Value = '  -2.00%  '
$ws.Cells.Item(12, 5).Value = '  -0.47%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '37.10'
$ws.Cells.Item(14, 5).Value = '  -0.28%  '
$ws.Cells.Item(15, 4).Value = '4.494.10'
$ws.Cells.Item(15, 5).Value = '  -1.66%  '
$ws.Cells.Item(16, 4).Value = '3.853.68'
$ws.Cells.Item(16, 5).Value = '  -1.63%  '
$ws.Cells.Item(17, 4).Value = '68.084.01'
$ws.Cells.Item(17, 5).Value = '  -1.26%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '7.51'
$ws.Cells.Item(18, 5).Value = '  +0.32%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '18.13'
$ws.Cells.Item(19, 5).Value = '  +5.98%  '
$ws.Cells.Item(20, 5).Value = '  -1.49%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.80'
$ws.Cells.Item(21, 5).Value = '  -3.19%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '470.03'
$ws.Cells.Item(22, 5).Value = '  -3.68%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.733'
$ws.Cells.Item(23, 5).Value = '  +1.23%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.0000161'
$ws.Cells.Item(24, 5).Value = '  -2.37%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '84.17'
$ws.Cells.Item(25, 5).Value = '  -0.28%  '
$ws.Cells.Item(26, 5).Value = '  -1.80%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.27'
$ws.Cells.Item(27, 5).Value = '  +1.18%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.01'
$ws.Cells.Item(28, 5).Value = '  -1.03%  '
$ws.Cells.Item(29, 5).Value = '  +0.04%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.92'
$ws.Cells.Item(30, 5).Value = '  -1.06%  '
$ws.Cells.Item(31, 4).Value = '3.998.80'
$ws.Cells.Item(31, 5).Value = '  -1.64%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.73'
$ws.Cells.Item(32, 5).Value = '  -1.86%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.31'
$ws.Cells.Item(33, 5).Value = '  -3.61%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '31.01'
$ws.Cells.Item(34, 5).Value = '  -3.96%  '
$ws.Cells.Item(35, 4).Value = '3.822.13'
$ws.Cells.Item(35, 5).Value = '  -0.97%  '
$ws.Cells.Item(36, 5).Value = '  -2.22%  '
$ws.Cells.Item(37, 5).Value = '  -1.07%  '
$ws.Cells.Item(38, 2).Value = 'Mantle'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.01'
$ws.Cells.Item(38, 5).Value = '  -2.59%  '
$ws.Cells.Item(39, 2).Value = 'Filecoin'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.93'
$ws.Cells.Item(39, 5).Value = '  +0.22%  '
$ws.Cells.Item(40, 2).Value = 'dogwifhat'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.29'
$ws.Cells.Item(40, 5).Value = '  +9.16%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.00'
$ws.Cells.Item(41, 5).Value = '  +0.07%  '
$ws.Cells.Item(42, 5).Value = '  -2.65%  '
$ws.Cells.Item(43, 5).Value = '  +0.12%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '427.08'
$ws.Cells.Item(44, 5).Value = '  -3.38%  '
$ws.Cells.Item(45, 5).Value = '  -0.01%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '47.35'
$ws.Cells.Item(46, 5).Value = '  -2.27%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '8.60'
$ws.Cells.Item(47, 5).Value = '  +1.02%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '143.20'
$ws.Cells.Item(48, 5).Value = '  +0.87%  '
$ws.Cells.Item(49, 5).Value = '  +15.32%  '
$ws.Cells.Item(50, 5).Value = '  +0.11%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '39.27'
$ws.Cells.Item(51, 5).Value = '  -0.02%  '
